$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.544.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.05%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.862.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.99%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.23%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'313.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.26%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.18%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4780"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.49%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3811"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.67%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07338"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.66%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9315"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.14%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.67%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07778"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.24%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.874.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.59%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.02%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.571"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.60%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'90.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.69%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.13%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008820"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.84%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.011"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'27.646.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.29%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.14%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.099"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.88%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.71%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.926"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.20%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'155.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +1.42%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.015"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.28%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'115.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.84%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'4.956"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.08866"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.03%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.332"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.50%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.7530"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.588"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.86%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.690"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.67%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.125"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +3.55%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.5622"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +7.34%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05337"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.981"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.25%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.032"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Aptos"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'8.509"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.85%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Algorand"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.1524"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.83%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +1.21%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4872"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.00%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -0.23%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'104.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.86%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.660"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.28%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'67.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.90%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +0.59%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.9096"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.10%  "
$ws.Range("E51").Style = "Normal"
